$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "Förändrad" (column C) from 2023-09-15 (45184) to 2023-09-17 (45186) ---
# Applies to every data row (2 through 344).
$lastRow = 344
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value() -ne $null) {
        $cell.Value = 45186
    }
}

# --- 2. Add a friendly display name (the "Beteckning" in column A) as the second
#        HYPERLINK() argument for every link cell in columns S..Y, rows 2..14 ---
$linkCols = 19, 20, 21, 22, 23, 24, 25   # S, T, U, V, W, X, Y

for ($r = 2; $r -le 14; $r++) {
    $label = $ws.Cells.Item($r, 1).Value()

    foreach ($c in $linkCols) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula) {
            $f = $cell.Formula()
            if ($f.Length -gt 0 -and $f.Substring($f.Length - 1) -eq ")" -and -not $f.Contains(",")) {
                $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $label + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
